$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @(2, 12, 0),
    @(3, 12, 0),
    @(4, 12, 0),
    @(5, 6, 0),
    @(6, 12, 0),
    @(7, 6, 0),
    @(8, 12, 0),
    @(9, 6, 0),
    @(10, 5, 0),
    @(11, 12, 0),
    @(12, 6, 0),
    @(13, 12, 0),
    @(14, 12, 0),
    @(15, 12, 0),
    @(16, 10, 20251119),
    @(17, 6, 0),
    @(18, 9, 0),
    @(19, 9, 0),
    @(20, 9, 0),
    @(21, 9, 0),
    @(22, 6, 0),
    @(23, 6, 0),
    @(24, 6, 0),
    @(25, 6, 0),
    @(26, 6, 0),
    @(27, 6, 0),
    @(28, 9, 0),
    @(29, 9, 0),
    @(30, 9, 0),
    @(31, 9, 0),
    @(32, 9, 0),
    @(33, 9, 0),
    @(34, 9, 0),
    @(35, 9, 0),
    @(37, 9, 0),
    @(38, 9, 0),
    @(39, 9, 0),
    @(40, 5, 0),
    @(41, 5, 0),
    @(42, 9, 0),
    @(43, 6, 0),
    @(44, 5, 0),
    @(45, 6, 0),
    @(46, 5, 0),
    @(47, 9, 0),
    @(48, 5, 0),
    @(49, 6, 0),
    @(50, 4, 0),
    @(51, 4, 0),
    @(52, 4, 0),
    @(53, 4, 0),
    @(54, 4, 0),
    @(55, 4, 0),
    @(56, 4, 0),
    @(57, 4, 0),
    @(58, 8, 0),
    @(59, 8, 0),
    @(60, 8, 0),
    @(61, 6, 0),
    @(62, 8, 0),
    @(63, 8, 0),
    @(64, 8, 0),
    @(65, 9, 0),
    @(66, 9, 0),
    @(67, 9, 0),
    @(68, 9, 0),
    @(69, 9, 0),
    @(70, 10, 20251119),
    @(71, 10, 20251119),
    @(72, 10, 20251119),
    @(73, 10, 20251119),
    @(74, 10, 20251119),
    @(75, 10, 20251119),
    @(76, 10, 20251119),
    @(77, 3, 0),
    @(78, 3, 0),
    @(79, 3, 0),
    @(80, 3, 0),
    @(81, 3, 0),
    @(82, 3, 0),
    @(83, 3, 0),
    @(84, 3, 0),
    @(85, 3, 0),
    @(86, 3, 0),
    @(87, 5, 0),
    @(88, 5, 0),
    @(89, 5, 0),
    @(90, 5, 0),
    @(91, 6, 0),
    @(92, 5, 0),
    @(93, 3, 0),
    @(94, 1, 0),
    @(95, 2, 0),
    @(96, 10, 20251119),
    @(97, 10, 20251119),
    @(98, 10, 20251119),
    @(99, 10, 20251119)
)

foreach ($change in $changes) {
    $row = $change[0]
    $newE = $change[1]
    $newF = $change[2]
    $ws.Cells.Item($row, 5).Value2 = $newE
    if ($newF -ne 0) {
        $ws.Cells.Item($row, 6).Value2 = $newF
    }
}
